# TC02_Canine_Filter_NeutStatus-No.xlsx
#
# The "CasesTab" query (row 2, column B of the "startup" sheet) is fixed:
# the trailing `, coalesce(co.cohort_description, '') AS `Cohort`` return
# column is dropped from the end of the Cypher query string (it matches
# the other "fixed" queries - Diagnosis / FileAssociation / FileFormat /
# FileType / NeuteredStatus / PrimeDiseaseSite - referenced in the commit
# message, all of which stop returning the extra Cohort column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesTabQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nWHERE demo.neutered_indicator IN [`"No`"]  `nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $casesTabQuery

# Rows re-flow (wrap-text autofit) once the query text is shorter; match
# the saved row heights / default row height from the re-saved workbook.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Selection left on B2 (the cell that was edited) instead of D3, and the
# view scrolled back up to the top of the sheet.
[void]$ws.Range("A1").Select()
[void]$ws.Range("B2").Select()
